# GitHub Actions cryptos-list refresh (Wed Apr 19 10:49:58 UTC 2023).
# Updates the per-coin Price (D) and Volume(1h) (E) snapshot columns on the
# 'Sheet1' crypto-ranking table, plus a reported rank swap between the
# WrappedEther / Chainlink rows (13 and 14 trade Coin name + Link).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 13/14: Coin name + Link swapped places (source re-ranked them) ---
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"

# --- Price column (D): these are free-text numeric strings (thousand-dot
#     grouping, fixed trailing zeros, etc.) that must round-trip exactly, so a
#     leading apostrophe forces text entry instead of Excel's numeric parsing.
#     ClearFormats() then drops the resulting quote-prefix marker so the cell
#     keeps the sheet's default (unstyled) formatting, matching the source. ---
$priceUpdates = @(
    @{ Cell = "D2"; Value = "29.354.01" }
    @{ Cell = "D3"; Value = "1.985.35" }
    @{ Cell = "D5"; Value = "329.92" }
    @{ Cell = "D7"; Value = "0.4943" }
    @{ Cell = "D8"; Value = "0.4196" }
    @{ Cell = "D9"; Value = "52.09" }
    @{ Cell = "D10"; Value = "0.08822" }
    @{ Cell = "D12"; Value = "23.23" }
    @{ Cell = "D13"; Value = "8.016" }
    @{ Cell = "D14"; Value = "1.965.88" }
    @{ Cell = "D15"; Value = "6.483" }
    @{ Cell = "D16"; Value = "96.03" }
    @{ Cell = "D18"; Value = "0.00001103" }
    @{ Cell = "D20"; Value = "19.68" }
    @{ Cell = "D21"; Value = "1.006" }
    @{ Cell = "D22"; Value = "5.940" }
    @{ Cell = "D23"; Value = "29.395.81" }
    @{ Cell = "D24"; Value = "11.85" }
    @{ Cell = "D25"; Value = "2.294" }
    @{ Cell = "D27"; Value = "20.50" }
    @{ Cell = "D28"; Value = "6.500" }
    @{ Cell = "D29"; Value = "2.337" }
    @{ Cell = "D30"; Value = "127.53" }
    @{ Cell = "D32"; Value = "0.09922" }
    @{ Cell = "D33"; Value = "1.562" }
    @{ Cell = "D34"; Value = "5.837" }
    @{ Cell = "D35"; Value = "3.782" }
    @{ Cell = "D36"; Value = "9.575" }
    @{ Cell = "D37"; Value = "0.02443" }
    @{ Cell = "D38"; Value = "0.06332" }
    @{ Cell = "D39"; Value = "1.280" }
    @{ Cell = "D41"; Value = "0.6481" }
    @{ Cell = "D42"; Value = "0.2063" }
    @{ Cell = "D46"; Value = "2.206" }
    @{ Cell = "D47"; Value = "1.261" }
    @{ Cell = "D48"; Value = "3.532" }
    @{ Cell = "D49"; Value = "0.00000000332" }
    @{ Cell = "D50"; Value = "0.06977" }
    @{ Cell = "D51"; Value = "1.141" }
)
foreach ($u in $priceUpdates) {
    $ws.Range($u.Cell).Value = "'" + $u.Value
    $ws.Range($u.Cell).ClearFormats()
}

# --- Volume(1h) column (E): padded '  +/-X.XX%  ' strings never parse as
#     numbers, so a plain text assignment round-trips exactly as-is. ---
$volumeUpdates = @(
    @{ Cell = "E2"; Value = "  -2.42%  " }
    @{ Cell = "E3"; Value = "  -6.30%  " }
    @{ Cell = "E4"; Value = "  -0.03%  " }
    @{ Cell = "E5"; Value = "  -4.82%  " }
    @{ Cell = "E6"; Value = "  +0.00%  " }
    @{ Cell = "E7"; Value = "  -4.85%  " }
    @{ Cell = "E8"; Value = "  -6.27%  " }
    @{ Cell = "E9"; Value = "  -4.05%  " }
    @{ Cell = "E10"; Value = "  -5.83%  " }
    @{ Cell = "E11"; Value = "  -5.55%  " }
    @{ Cell = "E12"; Value = "  -8.64%  " }
    @{ Cell = "E13"; Value = "  -7.80%  " }
    @{ Cell = "E14"; Value = "  -6.87%  " }
    @{ Cell = "E15"; Value = "  -7.15%  " }
    @{ Cell = "E16"; Value = "  -6.39%  " }
    @{ Cell = "E17"; Value = "  +0.07%  " }
    @{ Cell = "E18"; Value = "  -5.80%  " }
    @{ Cell = "E19"; Value = "  -0.89%  " }
    @{ Cell = "E20"; Value = "  -8.95%  " }
    @{ Cell = "E21"; Value = "  +0.02%  " }
    @{ Cell = "E22"; Value = "  -5.79%  " }
    @{ Cell = "E23"; Value = "  -2.39%  " }
    @{ Cell = "E24"; Value = "  -7.08%  " }
    @{ Cell = "E25"; Value = "  -1.70%  " }
    @{ Cell = "E26"; Value = "  -3.36%  " }
    @{ Cell = "E28"; Value = "  -2.88%  " }
    @{ Cell = "E29"; Value = "  -8.13%  " }
    @{ Cell = "E30"; Value = "  -5.00%  " }
    @{ Cell = "E31"; Value = "  -9.43%  " }
    @{ Cell = "E32"; Value = "  -6.14%  " }
    @{ Cell = "E33"; Value = "  -12.93%  " }
    @{ Cell = "E34"; Value = "  -7.14%  " }
    @{ Cell = "E35"; Value = "  -4.77%  " }
    @{ Cell = "E36"; Value = "  -10.88%  " }
    @{ Cell = "E37"; Value = "  -7.45%  " }
    @{ Cell = "E38"; Value = "  -7.94%  " }
    @{ Cell = "E39"; Value = "  -3.73%  " }
    @{ Cell = "E40"; Value = "  -7.85%  " }
    @{ Cell = "E41"; Value = "  -9.21%  " }
    @{ Cell = "E42"; Value = "  -8.51%  " }
    @{ Cell = "E43"; Value = "  +0.06%  " }
    @{ Cell = "E44"; Value = "  -8.69%  " }
    @{ Cell = "E46"; Value = "  -7.98%  " }
    @{ Cell = "E47"; Value = "  +0.11%  " }
    @{ Cell = "E48"; Value = "  -3.04%  " }
    @{ Cell = "E49"; Value = "  -6.78%  " }
    @{ Cell = "E50"; Value = "  -2.94%  " }
)
foreach ($u in $volumeUpdates) {
    $ws.Range($u.Cell).Value = $u.Value
}
